$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete Inflammatory-Mac / Resolving-Mac target-cluster rows (old rows 8-11)
$ws.Rows("8:11").Delete() | Out-Null

# Refresh remaining rows (2-7) with the re-run (new TPM) NATMI output
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Btc"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4875756666666667
$ws.Range("H2").Value = 1.462727
$ws.Range("I2").Value = 0.233114009085334
$ws.Range("J2").Value = 0.2331140090853341
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.232494
$ws.Range("N2").Value = 6.697482000000001
$ws.Range("O2").Value = 0.2983035867032896
$ws.Range("P2").Value = 0.2983035867032895
$ws.Range("Q2").Value = 1.088509750379334
$ws.Range("R2").Value = 9.796587753414002
$ws.Range("S2").Value = 0.06953874502093838
$ws.Range("T2").Value = 0.06953874502093838

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Btc"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4875756666666667
$ws.Range("H3").Value = 1.462727
$ws.Range("I3").Value = 0.233114009085334
$ws.Range("J3").Value = 0.2331140090853341
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.575715666666667
$ws.Range("N3").Value = 10.727147
$ws.Range("O3").Value = 0.4777835050834676
$ws.Range("P3").Value = 0.4777835050834675
$ws.Range("Q3").Value = 1.743431949985445
$ws.Range("R3").Value = 15.690887549869
$ws.Range("S3").Value = 0.1113780283448502
$ws.Range("T3").Value = 0.1113780283448502

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Btc"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4875756666666667
$ws.Range("H4").Value = 1.462727
$ws.Range("I4").Value = 0.233114009085334
$ws.Range("J4").Value = 0.2331140090853341
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.675756666666667
$ws.Range("N4").Value = 5.02727
$ws.Range("O4").Value = 0.2239129082132429
$ws.Range("P4").Value = 0.2239129082132428
$ws.Range("Q4").Value = 0.8170581739211111
$ws.Range("R4").Value = 7.35352356529
$ws.Range("S4").Value = 0.05219723571954546
$ws.Range("T4").Value = 0.05219723571954546

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Btc"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.604000333333333
$ws.Range("H5").Value = 4.812001
$ws.Range("I5").Value = 0.7668859909146659
$ws.Range("J5").Value = 0.7668859909146659
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.232494
$ws.Range("N5").Value = 6.697482000000001
$ws.Range("O5").Value = 0.2983035867032896
$ws.Range("P5").Value = 0.2983035867032895
$ws.Range("Q5").Value = 3.580921120164668
$ws.Range("R5").Value = 32.22829008148201
$ws.Range("S5").Value = 0.2287648416823512
$ws.Range("T5").Value = 0.2287648416823512

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Btc"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.604000333333333
$ws.Range("H6").Value = 4.812001
$ws.Range("I6").Value = 0.7668859909146659
$ws.Range("J6").Value = 0.7668859909146659
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.575715666666667
$ws.Range("N6").Value = 10.727147
$ws.Range("O6").Value = 0.4777835050834676
$ws.Range("P6").Value = 0.4777835050834675
$ws.Range("Q6").Value = 5.735449121238556
$ws.Range("R6").Value = 51.61904209114701
$ws.Range("S6").Value = 0.3664054767386173
$ws.Range("T6").Value = 0.3664054767386172

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Btc"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.604000333333333
$ws.Range("H7").Value = 4.812001
$ws.Range("I7").Value = 0.7668859909146659
$ws.Range("J7").Value = 0.7668859909146659
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.675756666666667
$ws.Range("N7").Value = 5.02727
$ws.Range("O7").Value = 0.2239129082132429
$ws.Range("P7").Value = 0.2239129082132428
$ws.Range("Q7").Value = 2.687914251918889
$ws.Range("R7").Value = 24.19122826727
$ws.Range("S7").Value = 0.1717156724936974
$ws.Range("T7").Value = 0.1717156724936974

